# Add 11 new day-columns (MR:NB) of mobility data to the "mobility" sheet,
# mirroring the style of the existing last column (MQ), and move the
# selection to the cell the author ended up leaving active (NH24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("MR","MS","MT","MU","MV","MW","MX","MY","MZ","NA","NB")

# Row 1 holds serial-date values formatted like the rest of the header
# row; row 2-5 hold the plain numeric data. Copying the existing MQ
# column into each new column first reproduces the same per-row cell
# style (date format on row 1, default/no style on rows 2-5) before we
# overwrite the copied values with the real ones below.
foreach ($col in $cols) {
    $ws.Range("MQ1:MQ5").Copy($ws.Range($col + "1"))
}

$row1 = @(44192, 44193, 44194, 44195, 44196, 44197, 44198, 44199, 44200, 44201, 44202)
$row2 = @(48.46, 61.02, 59.15, 51.9, 45.89, 43.48, 53.79, 46.34, 57.3, 63.1, 52.49)
$row3 = @(38.96, 48.54, 47.5, 47.95, 41.79, 31.83, 45.46, 35.35, 44.94, 44.17, 42.69)
$row4 = @(53.63, 66.7, 65.83, 64.94, 57.36, 50.2, 56.19, 54.99, 65.8, 67.52, 63.51)
$row5 = @(25.07, 28.77, 29.24, 28.61, 25.89, 22.29, 21.92, 24.82, 28.13, 30.37, 27.93)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    $ws.Range($col + "1").Value = $row1[$i]
    $ws.Range($col + "2").Value = $row2[$i]
    $ws.Range($col + "3").Value = $row3[$i]
    $ws.Range($col + "4").Value = $row4[$i]
    $ws.Range($col + "5").Value = $row5[$i]
}

# The author's sheet view ended up scrolled to the new columns with a
# single cell selected further to the right of the data.
[void]$ws.Range("NH24").Select()
$excel.ActiveWindow.ScrollColumn = 355
